# Fruta / hortaliza, semanal
# Inserts two new daily price records (rows 14 and 17) into the weekly
# "Vega Modelo de Temuco - Damasco" sheet, pushing the existing records
# down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows at their target positions.
# Inserting at 14 first, then at 17, lands both blank rows exactly where
# the new records belong (rows 15/16 become the former 14/15, and rows
# 18.. become the former 16..).
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(17).Insert()

# New record now sitting at row 14.
$row14 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44524, 9, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103003, 'Damasco', 'Castle Brite', 'Primera', 145, 18000, 20000, 19103, '$/bandeja 7 kilos', 'Provincia de San Felipe de Aconcagua', 2729, 7)
for ($i = 0; $i -lt $row14.Length; $i++) {
    $ws.Cells.Item(14, $i + 1).Value = $row14[$i]
}

# New record now sitting at row 17.
$row17 = @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44533, 9, 'Fruta', 100103, 'Frutos de hueso (carozo)', 100103003, 'Damasco', 'Castle Brite', 'Primera', 65, 1500, 1500, 1500, '$/bandeja 10 kilos', 'Provincia de San Felipe de Aconcagua', 150, 10)
for ($i = 0; $i -lt $row17.Length; $i++) {
    $ws.Cells.Item(17, $i + 1).Value = $row17[$i]
}

# Apply the same date display format used by the rest of column D to the
# date cells of the two freshly inserted rows.
$ws.Range("D14").NumberFormat = $ws.Range("D13").NumberFormat
$ws.Range("D17").NumberFormat = $ws.Range("D13").NumberFormat
